$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Cobertura de la Prueba")

# Fill in the missing week's actual UAT data (row 7: Planificadas/Disponibles/Ejecutadas/Ejecutadas OK)
$ws.Range("C7").Value = 20
$ws.Range("D7").Value = 19
$ws.Range("E7").Value = 19
$ws.Range("F7").Value = 19

# The "Sprint" markers in column A now align one row further down
$ws.Range("A7").Value = $null
$ws.Range("A8").Value = 3
$ws.Range("A9").Value = $null
$ws.Range("A10").Value = 4
$ws.Range("A11").Value = $null
$ws.Range("A12").Value = 5

# Row 10 moves out of the shaded band and matches the plain formatting used above it
$ws.Range("A10:F10").HorizontalAlignment = -4108

# Add the new week row at the bottom of the table (Semana 11)
$ws.Rows.Item(13).Insert()
$ws.Range("B13").Value = 11

# Grow the table to include the new row
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A2:F13"))

# Match the author's final selection
$ws.Range("F8").Select()

"done"
